$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, restricted to the used A:K columns only,
# so we don't touch/extend formatting across the full 16384-column row.
$ws.Range("A1:K1").Insert() | Out-Null

# Copy the style that used to belong to the (now shifted-down) header row,
# which is now on row 2, onto the brand-new row 1, keeping header look
# (bold font, centered/top alignment, thin border) on row 1.
$ws.Range("A2:K2").Copy() | Out-Null
$ws.Range("A1:K1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 1 becomes a simple numeric index: 0, 1, 2, ..., 10
for ($col = 1; $col -le 11; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# Row 2 (now holding what used to be the header labels) should look like a
# normal data row, not like the header: clear its special formatting.
$ws.Range("A2:K2").ClearFormats()

# In the original header row, H1 was blank, and J1/K1 held the internal
# field names "thread_size" / "material_surface"; once shifted to row 2
# those internal-only labels are removed, leaving H2, J2, K2 blank.
$ws.Cells.Item(2, 8).Value = $null
$ws.Cells.Item(2, 10).Value = $null
$ws.Cells.Item(2, 11).Value = $null
